$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the summary row (row 2) with the upgraded purchase analysis values
$ws.Range("A2").Value = 43969
$ws.Range("B2").Value = 43971
$ws.Range("C2").Value = 68659.80976688999
$ws.Range("D2").Value = 47785.18778951999
$ws.Range("E2").Value = 18835.08
$ws.Range("F2").Value = 2170.08
$ws.Range("G2").Value = 51694.66810475999
$ws.Range("H2").Value = 48689.86886209999
$ws.Range("I2").Value = 0.2470898436760764
